# Auto-generated edit script: updates currentAveragePrice / Leve profit
# columns (H:N) across all 8 job sheets to refreshed market-board values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2749
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H6").Value = 174.625
$ws.Range("I6").Value = 165.6
$ws.Range("J6").Value = 189.66667
$ws.Range("K6").Value = 496.8
$ws.Range("L6").Value = 569.00001
$ws.Range("M6").Value = -384.8
$ws.Range("N6").Value = -793.00001
$ws.Range("H9").Value = 440.6316
$ws.Range("I9").Value = 391.6
$ws.Range("J9").Value = 624.5
$ws.Range("K9").Value = 391.6
$ws.Range("L9").Value = 624.5
$ws.Range("M9").Value = -222.6
$ws.Range("N9").Value = -962.5
$ws.Range("H11").Value = 90909670
$ws.Range("I11").Value = 90909670
$ws.Range("K11").Value = 90909670
$ws.Range("M11").Value = -90909530
$ws.Range("H12").Value = 617.8570999999999
$ws.Range("I12").Value = 381.2857
$ws.Range("J12").Value = 854.4286
$ws.Range("K12").Value = 381.2857
$ws.Range("L12").Value = 854.4286
$ws.Range("M12").Value = -211.2857
$ws.Range("N12").Value = -1194.4286
$ws.Range("H15").Value = 1182.1864
$ws.Range("I15").Value = 1182.1864
$ws.Range("K15").Value = 3546.5592
$ws.Range("M15").Value = -3377.5592
$ws.Range("H21").Value = 35000
$ws.Range("J21").Value = 35000
$ws.Range("L21").Value = 35000
$ws.Range("N21").Value = -35936
$ws.Range("H23").Value = 35000
$ws.Range("J23").Value = 35000
$ws.Range("L23").Value = 35000
$ws.Range("N23").Value = -35468
$ws.Range("H29").Value = 303.75
$ws.Range("J29").Value = 400
$ws.Range("L29").Value = 1200
$ws.Range("N29").Value = -1762
$ws.Range("H33").Value = 422
$ws.Range("J33").Value = 476.2
$ws.Range("L33").Value = 476.2
$ws.Range("N33").Value = -934.2
$ws.Range("H38").Value = 1713.5
$ws.Range("I38").Value = 267.125
$ws.Range("J38").Value = 7499
$ws.Range("K38").Value = 801.375
$ws.Range("L38").Value = 22497
$ws.Range("M38").Value = -429.375
$ws.Range("N38").Value = -23241
$ws.Range("H43").Value = 5880.273
$ws.Range("I43").Value = 6098.7144
$ws.Range("J43").Value = 5498
$ws.Range("K43").Value = 6098.7144
$ws.Range("L43").Value = 5498
$ws.Range("M43").Value = -6029.7144
$ws.Range("N43").Value = -5636
$ws.Range("H58").Value = 4895.625
$ws.Range("I58").Value = 1055
$ws.Range("J58").Value = 7200
$ws.Range("K58").Value = 3165
$ws.Range("L58").Value = 21600
$ws.Range("M58").Value = -3015
$ws.Range("N58").Value = -21900
$ws.Range("H62").Value = 3781.8572
$ws.Range("I62").Value = 3412.1667
$ws.Range("K62").Value = 3412.1667
$ws.Range("M62").Value = -2788.1667
$ws.Range("H65").Value = 3781.8572
$ws.Range("I65").Value = 3412.1667
$ws.Range("K65").Value = 17060.8335
$ws.Range("M65").Value = -13940.8335
$ws.Range("H69").Value = 13074.5
$ws.Range("J69").Value = 13666.267
$ws.Range("L69").Value = 40998.801
$ws.Range("N69").Value = -42746.801
$ws.Range("H72").Value = 13074.5
$ws.Range("J72").Value = 13666.267
$ws.Range("L72").Value = 122996.403
$ws.Range("N72").Value = -131732.403
$ws.Range("H76").Value = 22249
$ws.Range("I76").Value = 22249
$ws.Range("K76").Value = 22249
$ws.Range("M76").Value = -21934
$ws.Range("H79").Value = 22249
$ws.Range("I79").Value = 22249
$ws.Range("K79").Value = 22249
$ws.Range("M79").Value = -21157
$ws.Range("H96").Value = 691372.4399999999
$ws.Range("I96").Value = 751.63635
$ws.Range("J96").Value = 1451055.2
$ws.Range("K96").Value = 2254.90905
$ws.Range("L96").Value = 4353165.6
$ws.Range("M96").Value = -881.9090500000002
$ws.Range("N96").Value = -4355911.6
$ws.Range("H103").Value = 45456412
$ws.Range("I103").Value = 1498.5
$ws.Range("K103").Value = 4495.5
$ws.Range("M103").Value = -3909.5
$ws.Range("H112").Value = 3728.8948
$ws.Range("J112").Value = 3974.303
$ws.Range("L112").Value = 11922.909
$ws.Range("N112").Value = -14138.909
$ws.Range("H118").Value = 3557.625
$ws.Range("I118").Value = 3610.7144
$ws.Range("J118").Value = 3186
$ws.Range("K118").Value = 10832.1432
$ws.Range("L118").Value = 9558
$ws.Range("M118").Value = -9175.143199999999
$ws.Range("N118").Value = -12872
$ws.Range("H125").Value = 1262.4286
$ws.Range("I125").Value = 543.75
$ws.Range("K125").Value = 4893.75
$ws.Range("M125").Value = -2433.75
$ws.Range("H132").Value = 1905.4186
$ws.Range("I132").Value = 1972.2821
$ws.Range("K132").Value = 5916.846299999999
$ws.Range("M132").Value = -3386.846299999999
$ws.Range("H135").Value = 1803.8
$ws.Range("I135").Value = 831.5172
$ws.Range("J135").Value = 30000
$ws.Range("K135").Value = 7483.6548
$ws.Range("L135").Value = 270000
$ws.Range("M135").Value = -4948.6548
$ws.Range("N135").Value = -275070
$ws.Range("H137").Value = 26320670
$ws.Range("I137").Value = 38465080
$ws.Range("K137").Value = 115395240
$ws.Range("M137").Value = -115392690
$ws.Range("H138").Value = 4486.6343
$ws.Range("I138").Value = 2438.5264
$ws.Range("J138").Value = 6255.4546
$ws.Range("K138").Value = 7315.5792
$ws.Range("L138").Value = 18766.3638
$ws.Range("M138").Value = -2175.5792
$ws.Range("N138").Value = -29046.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1469
$ws.Range("I2").Value = 1331.1428
$ws.Range("J2").Value = 1606.8572
$ws.Range("K2").Value = 1331.1428
$ws.Range("L2").Value = 1606.8572
$ws.Range("M2").Value = -1218.1428
$ws.Range("N2").Value = -1832.8572
$ws.Range("H6").Value = 500002100
$ws.Range("I6").Value = 4200
$ws.Range("J6").Value = 1000000000
$ws.Range("K6").Value = 4200
$ws.Range("L6").Value = 1000000000
$ws.Range("M6").Value = -4027
$ws.Range("N6").Value = -1000000346
$ws.Range("H13").Value = 333335330
$ws.Range("J13").Value = 500002000
$ws.Range("L13").Value = 500002000
$ws.Range("N13").Value = -500002288
$ws.Range("H28").Value = 8197.857
$ws.Range("I28").Value = 8197.857
$ws.Range("K28").Value = 8197.857
$ws.Range("M28").Value = -8005.857
$ws.Range("H32").Value = 10956.37
$ws.Range("I32").Value = 10605.177
$ws.Range("K32").Value = 10605.177
$ws.Range("M32").Value = -10318.177
$ws.Range("H36").Value = 502022500
$ws.Range("I36").Value = 504000000
$ws.Range("J36").Value = 500045000
$ws.Range("K36").Value = 504000000
$ws.Range("L36").Value = 500045000
$ws.Range("M36").Value = -503999654
$ws.Range("N36").Value = -500045692
$ws.Range("H45").Value = 5099.3335
$ws.Range("I45").Value = 2200
$ws.Range("J45").Value = 5679.2
$ws.Range("K45").Value = 2200
$ws.Range("L45").Value = 5679.2
$ws.Range("M45").Value = -1823
$ws.Range("N45").Value = -6433.2
$ws.Range("H61").Value = 22194310
$ws.Range("I61").Value = 25932742
$ws.Range("K61").Value = 25932742
$ws.Range("M61").Value = -25932530
$ws.Range("H74").Value = 53982.53
$ws.Range("I74").Value = 57200.188
$ws.Range("K74").Value = 57200.188
$ws.Range("M74").Value = -56326.188
$ws.Range("H77").Value = 53982.53
$ws.Range("I77").Value = 57200.188
$ws.Range("K77").Value = 286000.94
$ws.Range("M77").Value = -281632.94
$ws.Range("H99").Value = 8197.857
$ws.Range("I99").Value = 8197.857
$ws.Range("K99").Value = 8197.857
$ws.Range("M99").Value = -5202.857
$ws.Range("H110").Value = 4860.364
$ws.Range("I110").Value = 4398.353
$ws.Range("K110").Value = 4398.353
$ws.Range("M110").Value = -2353.353
$ws.Range("H116").Value = 1469
$ws.Range("I116").Value = 1331.1428
$ws.Range("J116").Value = 1606.8572
$ws.Range("K116").Value = 1331.1428
$ws.Range("L116").Value = 1606.8572
$ws.Range("M116").Value = 962.8571999999999
$ws.Range("N116").Value = -6194.8572
$ws.Range("H122").Value = 3095.3
$ws.Range("I122").Value = 2669.5833
$ws.Range("K122").Value = 8008.749899999999
$ws.Range("M122").Value = -5558.749899999999
$ws.Range("H132").Value = 4003108.5
$ws.Range("I132").Value = 3237.9583
$ws.Range("J132").Value = 100000000
$ws.Range("K132").Value = 9713.874899999999
$ws.Range("L132").Value = 300000000
$ws.Range("M132").Value = -7183.874899999999
$ws.Range("N132").Value = -300005060
$ws.Range("H136").Value = 22194310
$ws.Range("I136").Value = 25932742
$ws.Range("K136").Value = 77798226
$ws.Range("M136").Value = -77795676

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1469
$ws.Range("I3").Value = 1331.1428
$ws.Range("J3").Value = 1606.8572
$ws.Range("K3").Value = 1331.1428
$ws.Range("L3").Value = 1606.8572
$ws.Range("M3").Value = -1217.1428
$ws.Range("N3").Value = -1834.8572
$ws.Range("H8").Value = 368
$ws.Range("I8").Value = 152
$ws.Range("J8").Value = 800
$ws.Range("K8").Value = 152
$ws.Range("L8").Value = 800
$ws.Range("M8").Value = -12
$ws.Range("N8").Value = -1080
$ws.Range("H86").Value = 43789.062
$ws.Range("I86").Value = 52425
$ws.Range("J86").Value = 6366.6665
$ws.Range("K86").Value = 52425
$ws.Range("L86").Value = 6366.6665
$ws.Range("M86").Value = -51302
$ws.Range("N86").Value = -8612.666499999999
$ws.Range("H89").Value = 43789.062
$ws.Range("I89").Value = 52425
$ws.Range("J89").Value = 6366.6665
$ws.Range("K89").Value = 262125
$ws.Range("L89").Value = 31833.3325
$ws.Range("M89").Value = -256509
$ws.Range("N89").Value = -43065.3325
$ws.Range("H92").Value = 59000
$ws.Range("J92").Value = 59000
$ws.Range("L92").Value = 59000
$ws.Range("N92").Value = -63992
$ws.Range("H134").Value = 4349282
$ws.Range("I134").Value = 1522.2727
$ws.Range("J134").Value = 100000000
$ws.Range("K134").Value = 4566.8181
$ws.Range("L134").Value = 300000000
$ws.Range("M134").Value = -2031.8181
$ws.Range("N134").Value = -300005070

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H9").Value = 221904.8
$ws.Range("J9").Value = 221904.8
$ws.Range("L9").Value = 221904.8
$ws.Range("N9").Value = -222240.8
$ws.Range("H31").Value = 27780816
$ws.Range("I31").Value = 31252546
$ws.Range("J31").Value = 6982.75
$ws.Range("K31").Value = 31252546
$ws.Range("L31").Value = 6982.75
$ws.Range("M31").Value = -31252251
$ws.Range("N31").Value = -7572.75
$ws.Range("H34").Value = 27780816
$ws.Range("I34").Value = 31252546
$ws.Range("J34").Value = 6982.75
$ws.Range("K34").Value = 31252546
$ws.Range("L34").Value = 6982.75
$ws.Range("M34").Value = -31252344
$ws.Range("N34").Value = -7386.75
$ws.Range("H52").Value = 75650
$ws.Range("J52").Value = 98475
$ws.Range("L52").Value = 98475
$ws.Range("N52").Value = -99063
$ws.Range("H99").Value = 15710.392
$ws.Range("I99").Value = 9505.4
$ws.Range("J99").Value = 20483.46
$ws.Range("K99").Value = 9505.4
$ws.Range("L99").Value = 20483.46
$ws.Range("M99").Value = -8007.4
$ws.Range("N99").Value = -23479.46
$ws.Range("H103").Value = 27196.727
$ws.Range("I103").Value = 12862.667
$ws.Range("J103").Value = 44397.6
$ws.Range("K103").Value = 12862.667
$ws.Range("L103").Value = 44397.6
$ws.Range("M103").Value = -11690.667
$ws.Range("N103").Value = -46741.6
$ws.Range("H126").Value = 15710.392
$ws.Range("I126").Value = 9505.4
$ws.Range("J126").Value = 20483.46
$ws.Range("K126").Value = 28516.2
$ws.Range("L126").Value = 61450.38
$ws.Range("M126").Value = -26046.2
$ws.Range("N126").Value = -66390.38
$ws.Range("H132").Value = 1142.25
$ws.Range("I132").Value = 999.6667
$ws.Range("J132").Value = 1398.9
$ws.Range("K132").Value = 2999.0001
$ws.Range("L132").Value = 4196.700000000001
$ws.Range("M132").Value = -469.0001000000002
$ws.Range("N132").Value = -9256.700000000001
$ws.Range("H134").Value = 1894.8
$ws.Range("I134").Value = 1693.8387
$ws.Range("J134").Value = 3452.25
$ws.Range("K134").Value = 5081.5161
$ws.Range("L134").Value = 10356.75
$ws.Range("M134").Value = -2546.5161
$ws.Range("N134").Value = -15426.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 333333700
$ws.Range("I7").Value = 1000000000
$ws.Range("J7").Value = 539.5
$ws.Range("K7").Value = 3000000000
$ws.Range("L7").Value = 1618.5
$ws.Range("M7").Value = -2999999888
$ws.Range("N7").Value = -1842.5
$ws.Range("H12").Value = 45457740
$ws.Range("I12").Value = 125003410
$ws.Range("J12").Value = 3075.1428
$ws.Range("K12").Value = 375010230
$ws.Range("L12").Value = 9225.428400000001
$ws.Range("M12").Value = -375010057
$ws.Range("N12").Value = -9571.428400000001
$ws.Range("H26").Value = 194.75
$ws.Range("I26").Value = 78.59999999999999
$ws.Range("J26").Value = 388.33334
$ws.Range("K26").Value = 235.8
$ws.Range("L26").Value = 1165.00002
$ws.Range("M26").Value = 52.20000000000002
$ws.Range("N26").Value = -1741.00002
$ws.Range("H80").Value = 23819048
$ws.Range("I80").Value = 27788222
$ws.Range("K80").Value = 83364666
$ws.Range("M80").Value = -83363730
$ws.Range("H81").Value = 4866.5
$ws.Range("I81").Value = 1733
$ws.Range("J81").Value = 8000
$ws.Range("K81").Value = 5199
$ws.Range("L81").Value = 24000
$ws.Range("M81").Value = -4076
$ws.Range("N81").Value = -26246
$ws.Range("H83").Value = 23819048
$ws.Range("I83").Value = 27788222
$ws.Range("K83").Value = 250093998
$ws.Range("M83").Value = -250089318
$ws.Range("H84").Value = 4866.5
$ws.Range("I84").Value = 1733
$ws.Range("J84").Value = 8000
$ws.Range("K84").Value = 15597
$ws.Range("L84").Value = 72000
$ws.Range("M84").Value = -9981
$ws.Range("N84").Value = -83232
$ws.Range("H92").Value = 212
$ws.Range("J92").Value = 199.33333
$ws.Range("L92").Value = 597.99999
$ws.Range("N92").Value = -3093.99999
$ws.Range("H94").Value = 20000
$ws.Range("J94").Value = 20000
$ws.Range("L94").Value = 60000
$ws.Range("N94").Value = -61352
$ws.Range("H129").Value = 2319.3928
$ws.Range("I129").Value = 2014.8182
$ws.Range("K129").Value = 6044.4546
$ws.Range("M129").Value = -1044.4546
$ws.Range("H131").Value = 2400.1365
$ws.Range("I131").Value = 3000
$ws.Range("K131").Value = 9000
$ws.Range("M131").Value = -3960
$ws.Range("H138").Value = 10299.208
$ws.Range("I138").Value = 12066.429
$ws.Range("J138").Value = 7825.1
$ws.Range("K138").Value = 36199.287
$ws.Range("L138").Value = 23475.3
$ws.Range("M138").Value = -31059.287
$ws.Range("N138").Value = -33755.3

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 20718.092
$ws.Range("I5").Value = 42966.332
$ws.Range("J5").Value = 12375
$ws.Range("K5").Value = 42966.332
$ws.Range("L5").Value = 12375
$ws.Range("M5").Value = -42854.332
$ws.Range("N5").Value = -12599
$ws.Range("H21").Value = 66684664
$ws.Range("I21").Value = 66684664
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 66684664
$ws.Range("L21").Value = 0
$ws.Range("M21").Value = -66684491
$ws.Range("N21").ClearContents()
$ws.Range("H30").Value = 66684664
$ws.Range("I30").Value = 66684664
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 66684664
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -66684559
$ws.Range("N30").ClearContents()
$ws.Range("H62").Value = 24000
$ws.Range("J62").Value = 24000
$ws.Range("L62").Value = 24000
$ws.Range("N62").Value = -25372
$ws.Range("H65").Value = 24000
$ws.Range("J65").Value = 24000
$ws.Range("L65").Value = 72000
$ws.Range("N65").Value = -78864
$ws.Range("H102").Value = 1224.4828
$ws.Range("I102").Value = 1135.6154
$ws.Range("J102").Value = 1994.6666
$ws.Range("K102").Value = 1135.6154
$ws.Range("L102").Value = 1994.6666
$ws.Range("M102").Value = 486.3846000000001
$ws.Range("N102").Value = -5238.6666
$ws.Range("H113").Value = 1853972.1
$ws.Range("I113").Value = 2341.8572
$ws.Range("J113").Value = 6174442.5
$ws.Range("K113").Value = 2341.8572
$ws.Range("L113").Value = 6174442.5
$ws.Range("M113").Value = -171.8571999999999
$ws.Range("N113").Value = -6178782.5
$ws.Range("H126").Value = 13044831
$ws.Range("I126").Value = 14178773
$ws.Range("K126").Value = 42536319
$ws.Range("M126").Value = -42533849
$ws.Range("H132").Value = 5970643
$ws.Range("I132").Value = 4964.731
$ws.Range("J132").Value = 31821914
$ws.Range("K132").Value = 14894.193
$ws.Range("L132").Value = 95465742
$ws.Range("M132").Value = -12364.193
$ws.Range("N132").Value = -95470802
$ws.Range("H136").Value = 17500
$ws.Range("J136").Value = 17500
$ws.Range("L136").Value = 52500
$ws.Range("N136").Value = -57600

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 10872
$ws.Range("I7").Value = 9979.235000000001
$ws.Range("J7").Value = 14666.25
$ws.Range("K7").Value = 9979.235000000001
$ws.Range("L7").Value = 14666.25
$ws.Range("M7").Value = -9867.235000000001
$ws.Range("N7").Value = -14890.25
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").ClearContents()
$ws.Range("H46").Value = 1195.3125
$ws.Range("I46").Value = 989.25
$ws.Range("J46").Value = 1401.375
$ws.Range("K46").Value = 989.25
$ws.Range("L46").Value = 1401.375
$ws.Range("M46").Value = -801.25
$ws.Range("N46").Value = -1777.375
$ws.Range("H61").Value = 3285.6667
$ws.Range("I61").Value = 2166.4167
$ws.Range("J61").Value = 7762.6665
$ws.Range("K61").Value = 2166.4167
$ws.Range("L61").Value = 7762.6665
$ws.Range("M61").Value = -1964.4167
$ws.Range("N61").Value = -8166.6665
$ws.Range("H68").Value = 3478555.2
$ws.Range("I68").Value = 5213083
$ws.Range("J68").Value = 9500
$ws.Range("K68").Value = 5213083
$ws.Range("L68").Value = 9500
$ws.Range("M68").Value = -5212334
$ws.Range("N68").Value = -10998
$ws.Range("H69").Value = 51880
$ws.Range("J69").Value = 51880
$ws.Range("L69").Value = 51880
$ws.Range("N69").Value = -53502
$ws.Range("H71").Value = 3478555.2
$ws.Range("I71").Value = 5213083
$ws.Range("J71").Value = 9500
$ws.Range("K71").Value = 26065415
$ws.Range("L71").Value = 47500
$ws.Range("M71").Value = -26061671
$ws.Range("N71").Value = -54988
$ws.Range("H72").Value = 51880
$ws.Range("J72").Value = 51880
$ws.Range("L72").Value = 155640
$ws.Range("N72").Value = -163752
$ws.Range("H100").Value = 41703900
$ws.Range("J100").Value = 62505100
$ws.Range("L100").Value = 62505100
$ws.Range("N100").Value = -62506182
$ws.Range("H113").Value = 3285.6667
$ws.Range("I113").Value = 2166.4167
$ws.Range("J113").Value = 7762.6665
$ws.Range("K113").Value = 2166.4167
$ws.Range("L113").Value = 7762.6665
$ws.Range("M113").Value = 3.583299999999781
$ws.Range("N113").Value = -12102.6665
$ws.Range("H119").Value = 89750
$ws.Range("J119").Value = 89750
$ws.Range("L119").Value = 89750
$ws.Range("N119").Value = -99426
$ws.Range("H122").Value = 3348.362
$ws.Range("I122").Value = 3296.3584
$ws.Range("K122").Value = 9889.075199999999
$ws.Range("M122").Value = -7439.075199999999
$ws.Range("H126").Value = 10872
$ws.Range("I126").Value = 9979.235000000001
$ws.Range("J126").Value = 14666.25
$ws.Range("K126").Value = 29937.705
$ws.Range("L126").Value = 43998.75
$ws.Range("M126").Value = -27467.705
$ws.Range("N126").Value = -48938.75
$ws.Range("H132").Value = 3259.25
$ws.Range("I132").Value = 2718.65
$ws.Range("J132").Value = 5962.25
$ws.Range("K132").Value = 8155.950000000001
$ws.Range("L132").Value = 17886.75
$ws.Range("M132").Value = -5625.950000000001
$ws.Range("N132").Value = -22946.75
$ws.Range("H136").Value = 3924.2
$ws.Range("I136").Value = 3693.5557
$ws.Range("J136").Value = 6000
$ws.Range("K136").Value = 11080.6671
$ws.Range("L136").Value = 18000
$ws.Range("M136").Value = -8530.667099999999
$ws.Range("N136").Value = -23100

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8125.696
$ws.Range("J62").Value = 10581.667
$ws.Range("L62").Value = 10581.667
$ws.Range("N62").Value = -11829.667
$ws.Range("H65").Value = 8125.696
$ws.Range("J65").Value = 10581.667
$ws.Range("L65").Value = 52908.335
$ws.Range("N65").Value = -59148.335
$ws.Range("H81").Value = 2282.8
$ws.Range("I81").Value = 1724.7693
$ws.Range("K81").Value = 3449.5386
$ws.Range("M81").Value = -2388.5386
$ws.Range("H84").Value = 2282.8
$ws.Range("I84").Value = 1724.7693
$ws.Range("K84").Value = 17247.693
$ws.Range("M84").Value = -11943.693
$ws.Range("H105").Value = 38332
$ws.Range("J105").Value = 38332
$ws.Range("L105").Value = 38332
$ws.Range("N105").Value = -45320
$ws.Range("H113").Value = 809.8333
$ws.Range("I113").Value = 595.25
$ws.Range("K113").Value = 1785.75
$ws.Range("M113").Value = 384.25
$ws.Range("H136").Value = 299793.78
$ws.Range("I136").Value = 6377.7407
$ws.Range("J136").Value = 1431541.4
$ws.Range("K136").Value = 19133.2221
$ws.Range("L136").Value = 4294624.199999999
$ws.Range("M136").Value = -16583.2221
$ws.Range("N136").Value = -4299724.199999999
